# Update metrics for all model rows (rows 2-26) with new values from the
# retrained ensemble ("atualizado todo o treinamento para o novo lm").
# Every data row shares the same new metric values across columns B:Q.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    0.9999824846980448,        # B r2
    0.9990072066285755,        # C r2_sup
    0.9997097675380739,        # D r2_test
    0.9999454258523077,        # E r2_val
    0.9999494320524768,        # F r2_vt
    0.00001634976554226276,    # G mse
    0.0009267290336323482,     # H mse_sup
    0.00007402531320678512,    # I mse_test
    0.00001428024827027335,    # J mse_val
    0.00004415278073852923,    # K mse_vt
    0.0002547310818618769,     # L mape
    0.004043484331892824,      # M rmse
    1.000016814689877,         # N r2_adj
    0.004215623831067973,      # O rsd
    120.0425940012932,         # P aic
    179.767509419835           # Q bic
)

$firstRow = 2
$lastRow = 26
$firstCol = 2   # column B
$lastCol = 17   # column Q

for ($row = $firstRow; $row -le $lastRow; $row++) {
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $ws.Cells.Item($row, $col).Value = $newValues[$col - $firstCol]
    }
}

Write-Output "Updated metrics for rows $firstRow to $lastRow"
